# Updates the cryptocurrency price/volume table to reflect the latest
# scrape performed by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") stores every value as literal text (e.g. "46.824.02"),
# so force a Text number format before writing the new figures. Without
# this, Excel would silently reinterpret values such as "302.16" or
# "1.00" as numbers and corrupt/round them.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '46.778.90'
$ws.Range("E2").Value = '  +5.19%  '

# Row 3
$ws.Range("D3").Value = '2.304.36'
$ws.Range("E3").Value = '  +3.80%  '

# Row 4
$ws.Range("E4").Value = '  -0.70%  '

# Row 5
$ws.Range("D5").Value = '302.16'
$ws.Range("E5").Value = '  +1.42%  '

# Row 6
$ws.Range("D6").Value = '101.81'
$ws.Range("E6").Value = '  +12.72%  '

# Row 7
$ws.Range("D7").Value = '0.571'
$ws.Range("E7").Value = '  +0.96%  '

# Row 8
$ws.Range("E8").Value = '  -0.56%  '

# Row 9
$ws.Range("E9").Value = '  +8.54%  '

# Row 10
$ws.Range("D10").Value = '36.74'
$ws.Range("E10").Value = '  +10.80%  '

# Row 11
$ws.Range("D11").Value = '0.0805'
$ws.Range("E11").Value = '  +2.20%  '

# Row 12
$ws.Range("D12").Value = '7.38'
$ws.Range("E12").Value = '  +5.81%  '

# Row 13
$ws.Range("E13").Value = '  +0.69%  '

# Row 14
$ws.Range("D14").Value = '2.655.04'
$ws.Range("E14").Value = '  +3.81%  '

# Row 15
$ws.Range("D15").Value = '2.307.80'
$ws.Range("E15").Value = '  -26.80%  '

# Row 16
$ws.Range("D16").Value = '14.02'
$ws.Range("E16").Value = '  +3.67%  '

# Row 17
$ws.Range("D17").Value = '0.821'
$ws.Range("E17").Value = '  +4.15%  '

# Row 18
$ws.Range("D18").Value = '46.719.94'
$ws.Range("E18").Value = '  +5.70%  '

# Row 19
$ws.Range("D19").Value = '13.44'
$ws.Range("E19").Value = '  +20.66%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0948'
$ws.Range("E20").Value = '  +4.62%  '

# Row 21
$ws.Range("D21").Value = '6.13'
$ws.Range("E21").Value = '  +3.40%  '

# Row 22
$ws.Range("D22").Value = '66.89'
$ws.Range("E22").Value = '  +4.24%  '

# Row 23
$ws.Range("D23").Value = '248.82'
$ws.Range("E23").Value = '  +5.09%  '

# Row 24
$ws.Range("E24").Value = '  +5.73%  '

# Row 25
$ws.Range("E25").Value = '  +5.41%  '

# Row 26
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -1.07%  '

# Row 27
$ws.Range("D27").Value = '44.82'
$ws.Range("E27").Value = '  +16.42%  '

# Row 28
$ws.Range("D28").Value = '2.28'
$ws.Range("E28").Value = '  +1.28%  '

# Row 29
$ws.Range("D29").Value = '9.91'
$ws.Range("E29").Value = '  +5.53%  '

# Row 30
$ws.Range("D30").Value = '20.18'
$ws.Range("E30").Value = '  +3.24%  '

# Row 31
$ws.Range("D31").Value = '5.81'
$ws.Range("E31").Value = '  +7.50%  '

# Row 32
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value = '147.13'
$ws.Range("E32").Value = '  -1.01%  '

# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.0801'
$ws.Range("E33").Value = '  +6.64%  '

# Row 34
$ws.Range("E34").Value = '  +2.28%  '

# Row 35
$ws.Range("D35").Value = '3.14'
$ws.Range("E35").Value = '  +9.71%  '

# Row 36
$ws.Range("E36").Value = '  +8.88%  '

# Row 37
$ws.Range("E37").Value = '  +2.40%  '

# Row 38
$ws.Range("D38").Value = '1.81'
$ws.Range("E38").Value = '  +8.31%  '

# Row 39
$ws.Range("D39").Value = '15.89'
$ws.Range("E39").Value = '  +19.93%  '

# Row 40
$ws.Range("E40").Value = '  +13.85%  '

# Row 41
$ws.Range("E41").Value = '  +10.26%  '

# Row 42
$ws.Range("E42").Value = '  +2.79%  '

# Row 43
$ws.Range("E43").Value = '  -0.72%  '

# Row 44
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '1.867.07'
$ws.Range("E44").Value = '  +2.40%  '

# Row 45
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Value = '1.98'
$ws.Range("E45").Value = '  +10.59%  '

# Row 46
$ws.Range("D46").Value = '88.01'
$ws.Range("E46").Value = '  +18.88%  '

# Row 47
$ws.Range("E47").Value = '  +9.45%  '

# Row 48
$ws.Range("D48").Value = '74.67'
$ws.Range("E48").Value = '  +10.45%  '

# Row 49
$ws.Range("D49").Value = '4.93'
$ws.Range("E49").Value = '  +11.71%  '

# Row 50
$ws.Range("D50").Value = '97.58'
$ws.Range("E50").Value = '  +3.05%  '

# Row 51
$ws.Range("D51").Value = '8.07'
$ws.Range("E51").Value = '  +5.49%  '
